$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("res_type_Energy_source")

# Row 2 : Biomasse / apartment
$ws.Range("C2").Value = 320487
$ws.Range("D2").Value = 239.4219237582472
$ws.Range("E2").Value = 66.95905227218795

# Row 10 : Biomasse / house
$ws.Range("C10").Value = 3775624
$ws.Range("D10").Value = 274.568821693265
$ws.Range("E10").Value = 112.8278356336831
